$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.143.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.13%  "

# Row 3
$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.054.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.61%  "

# Row 4
$ws.Range("D4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.80%  "

# Row 6
$ws.Range("D6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.654"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.24%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +15.97%  "

# Row 9
$ws.Range("D9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.74"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.43%  "

# Row 10
$ws.Range("D10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.375"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.42%  "

# Row 11
$ws.Range("D11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0785"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.80%  "

# Row 12
$ws.Range("E12").Value = "  +5.34%  "

# Row 13
$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.97%  "

# Row 14
$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.349.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.51%  "

# Row 15
$ws.Range("D15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.813"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.56%  "

# Row 16
$ws.Range("D16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.64%  "

# Row 17
$ws.Range("D17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.049.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.45%  "

# Row 18
$ws.Range("D18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.100.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.99%  "

# Row 19
$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0905"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.01%  "

# Row 20
$ws.Range("D20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.02%  "

# Row 21
$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.49%  "

# Row 22
$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.90%  "

# Row 23
$ws.Range("D23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.07%  "

# Row 24
$ws.Range("E24").Value = "  -0.05%  "

# Row 25
$ws.Range("E25").Value = "  -2.75%  "

# Row 26
$ws.Range("D26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.13"
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("D27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.75%  "

# Row 28
$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.87%  "

# Row 29
$ws.Range("D29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.76%  "

# Row 30
$ws.Range("E30").Value = "  -0.38%  "

# Row 31
$ws.Range("E31").Value = "  +0.84%  "

# Row 32
$ws.Range("E32").Value = "  +12.02%  "

# Row 33
$ws.Range("D33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0624"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.94%  "

# Row 34
$ws.Range("D34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.50%  "

# Row 35
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0877"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.01%  "

# Row 36
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.11%  "

# Row 37
$ws.Range("E37").Value = "  -5.03%  "

# Row 38
$ws.Range("E38").Value = "  -8.46%  "

# Row 39
$ws.Range("D39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.57%  "

# Row 40
$ws.Range("B40").Value = "Cronos"
$ws.Range("C40").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.103"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +21.90%  "

# Row 41
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.63%  "

# Row 42
$ws.Range("B42").Value = "Gas"
$ws.Range("C42").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -43.46%  "

# Row 43
$ws.Range("E43").Value = "  -0.85%  "

# Row 44
$ws.Range("D44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.32%  "

# Row 45
$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.30%  "

# Row 46
$ws.Range("D46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.40%  "

# Row 47
$ws.Range("D47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +29.72%  "

# Row 48
$ws.Range("D48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.295.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.87%  "

# Row 49
$ws.Range("D49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.66%  "

# Row 50
$ws.Range("D50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.97%  "

# Row 51
$ws.Range("D51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.44%  "

